$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 345.26666
$ws.Range("I80").Value = 335.66666
$ws.Range("K80").Value = 1006.99998
$ws.Range("M80").Value = -8.999979999999937

$ws.Range("H83").Value = 345.26666
$ws.Range("I83").Value = 335.66666
$ws.Range("K83").Value = 3020.99994
$ws.Range("M83").Value = 1971.00006

$ws.Range("H111").Value = 2620.7273
$ws.Range("J111").Value = 2653
$ws.Range("L111").Value = 7959
$ws.Range("N111").Value = -14093

$ws.Range("H125").Value = 16581.125
$ws.Range("J125").Value = 499.66666
$ws.Range("L125").Value = 4496.99994
$ws.Range("N125").Value = -9416.99994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 2513
$ws.Range("I33").Value = 2513
$ws.Range("K33").Value = 2513
$ws.Range("M33").Value = -2184

$ws.Range("H45").Value = 4843.8887
$ws.Range("I45").Value = 2523.75
$ws.Range("J45").Value = 6700
$ws.Range("K45").Value = 2523.75
$ws.Range("L45").Value = 6700
$ws.Range("M45").Value = -2146.75
$ws.Range("N45").Value = -7454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1751.25
$ws.Range("I94").Value = 1463.9412
$ws.Range("J94").Value = 2195.2727
$ws.Range("K94").Value = 1463.9412
$ws.Range("L94").Value = 2195.2727
$ws.Range("M94").Value = -1012.9412
$ws.Range("N94").Value = -3097.2727

$ws.Range("H134").Value = 6859.077
$ws.Range("I134").Value = 7014.25
$ws.Range("J134").Value = 4997
$ws.Range("K134").Value = 21042.75
$ws.Range("L134").Value = 14991
$ws.Range("M134").Value = -18507.75
$ws.Range("N134").Value = -20061

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 649.44446
$ws.Range("I5").Value = 228.4
$ws.Range("J5").Value = 1175.75
$ws.Range("K5").Value = 228.4
$ws.Range("L5").Value = 1175.75
$ws.Range("M5").Value = -116.4
$ws.Range("N5").Value = -1399.75

$ws.Range("H17").Value = 215
$ws.Range("I17").Value = 215
$ws.Range("K17").Value = 215
$ws.Range("M17").Value = -41

$ws.Range("H58").Value = 10582.6
$ws.Range("I58").Value = 8131.3335
$ws.Range("K58").Value = 8131.3335
$ws.Range("M58").Value = -7928.3335

$ws.Range("H107").Value = 684.78125
$ws.Range("I107").Value = 635.3461
$ws.Range("K107").Value = 635.3461
$ws.Range("M107").Value = 1284.6539

$ws.Range("H134").Value = 6120.28
$ws.Range("I134").Value = 5399.65
$ws.Range("K134").Value = 16198.95
$ws.Range("M134").Value = -13663.95

$ws.Range("H136").Value = 10582.6
$ws.Range("I136").Value = 8131.3335
$ws.Range("K136").Value = 24394.0005
$ws.Range("M136").Value = -21844.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I34").Value = 233.33333
$ws.Range("J34").Value = 448.42856
$ws.Range("K34").Value = 699.99999
$ws.Range("L34").Value = 1345.28568
$ws.Range("M34").Value = -615.99999
$ws.Range("N34").Value = -1513.28568

$ws.Range("H128").Value = 666829.75
$ws.Range("I128").Value = 666829.75
$ws.Range("K128").Value = 2000489.25
$ws.Range("M128").Value = -1995509.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3687.889
$ws.Range("J122").Value = 2926.25
$ws.Range("L122").Value = 8778.75
$ws.Range("N122").Value = -13678.75

$ws.Range("H126").Value = 4714.875
$ws.Range("I126").Value = 4362.533
$ws.Range("K126").Value = 13087.599
$ws.Range("M126").Value = -10617.599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3843.5
$ws.Range("I7").Value = 3874.6667
$ws.Range("K7").Value = 3874.6667
$ws.Range("M7").Value = -3762.6667

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null

$ws.Range("H32").Value = 9666.333000000001
$ws.Range("I32").Value = 9666.333000000001
$ws.Range("K32").Value = 9666.333000000001
$ws.Range("M32").Value = -9349.333000000001

$ws.Range("H40").Value = 2243.5
$ws.Range("I40").Value = 1988
$ws.Range("J40").Value = 2499
$ws.Range("K40").Value = 1988
$ws.Range("L40").Value = 2499
$ws.Range("M40").Value = -1852
$ws.Range("N40").Value = -2771

$ws.Range("H55").Value = 666.94116
$ws.Range("I55").Value = 726.2727
$ws.Range("K55").Value = 726.2727
$ws.Range("M55").Value = -553.2727

$ws.Range("H61").Value = 113667
$ws.Range("I61").Value = 169067.33
$ws.Range("K61").Value = 169067.33
$ws.Range("M61").Value = -168865.33

$ws.Range("H113").Value = 113667
$ws.Range("I113").Value = 169067.33
$ws.Range("K113").Value = 169067.33
$ws.Range("M113").Value = -166897.33

$ws.Range("H122").Value = 4442
$ws.Range("I122").Value = 4432.3335
$ws.Range("J122").Value = 4449.25
$ws.Range("K122").Value = 13297.0005
$ws.Range("L122").Value = 13347.75
$ws.Range("M122").Value = -10847.0005
$ws.Range("N122").Value = -18247.75

$ws.Range("H126").Value = 3843.5
$ws.Range("I126").Value = 3874.6667
$ws.Range("K126").Value = 11624.0001
$ws.Range("M126").Value = -9154.000100000001

$ws.Range("H132").Value = 14761.8
$ws.Range("I132").Value = 16528.268
$ws.Range("K132").Value = 49584.804
$ws.Range("M132").Value = -47054.804

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2887

$ws.Range("H32").Value = 19200
$ws.Range("I32").Value = 8500
$ws.Range("J32").Value = 29900
$ws.Range("K32").Value = 8500
$ws.Range("L32").Value = 29900
$ws.Range("M32").Value = -8183
$ws.Range("N32").Value = -30534

$ws.Range("H34").Value = 15333
$ws.Range("I34").Value = 12999.5
$ws.Range("K34").Value = 12999.5
$ws.Range("M34").Value = -12796.5

$ws.Range("H51").Value = 49999
$ws.Range("J51").Value = 49999
$ws.Range("L51").Value = 49999
$ws.Range("N51").Value = -51019

$ws.Range("H74").Value = 20671.375
$ws.Range("I74").Value = 15949
$ws.Range("J74").Value = 21346
$ws.Range("K74").Value = 15949
$ws.Range("L74").Value = 21346
$ws.Range("M74").Value = -15013
$ws.Range("N74").Value = -23218

$ws.Range("H77").Value = 20671.375
$ws.Range("I77").Value = 15949
$ws.Range("J77").Value = 21346
$ws.Range("K77").Value = 47847
$ws.Range("L77").Value = 64038
$ws.Range("M77").Value = -43167
$ws.Range("N77").Value = -73398

$ws.Range("H113").Value = 787.3333
$ws.Range("J113").Value = 524.5
$ws.Range("L113").Value = 1573.5
$ws.Range("N113").Value = -5913.5

$ws.Range("H122").Value = 2464.3333
$ws.Range("I122").Value = 1577.2667
$ws.Range("K122").Value = 4731.800099999999
$ws.Range("M122").Value = -2281.800099999999

$ws.Range("H126").Value = 7005
$ws.Range("J126").Value = 7168
$ws.Range("L126").Value = 21504
$ws.Range("N126").Value = -26444

$ws.Range("H132").Value = 6013.926
$ws.Range("I132").Value = 5332.2383
$ws.Range("K132").Value = 15996.7149
$ws.Range("M132").Value = -13466.7149
